$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 174
$ws.Range("F6").Value = 2759
$ws.Range("F8").Value = 1647
$ws.Range("F9").Value = 7467
$ws.Range("F11").Value = 7661
$ws.Range("F13").Value = 36
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = 6191
$ws.Range("F16").Value = 3270
$ws.Range("F17").Value = 3634
$ws.Range("F19").Value = 14
$ws.Range("F21").Value = 31
$ws.Range("F24").Value = 287
$ws.Range("F25").Value = 289
$ws.Range("F26").Value = 3648
$ws.Range("F28").Value = 343
$ws.Range("F30").Value = 263
$ws.Range("F31").Value = 1103
$ws.Range("F35").Value = 1474
$ws.Range("F37").Value = 21
$ws.Range("F39").Value = 3296
$ws.Range("F41").Value = 246
$ws.Range("F45").Value = 1290
$ws.Range("F47").Value = 526
$ws.Range("F48").Value = 596

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 58
$ws.Range("F6").Value = 44

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 122

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 174
$ws.Range("F6").Value = 58
$ws.Range("F8").Value = 122
$ws.Range("F9").Value = 1647
$ws.Range("F11").Value = 44
$ws.Range("F12").Value = 7467
$ws.Range("F13").Value = 7661
$ws.Range("F15").Value = 6191
$ws.Range("F16").Value = 3270
$ws.Range("F17").Value = 3634
$ws.Range("F19").Value = 14
$ws.Range("F21").Value = 31
$ws.Range("F23").Value = 287
$ws.Range("F26").Value = 289
$ws.Range("F27").Value = 3648
$ws.Range("F30").Value = 343
$ws.Range("F32").Value = 263
$ws.Range("F36").Value = 1474
$ws.Range("F38").Value = 21
$ws.Range("F40").Value = 3296
$ws.Range("F42").Value = 246
$ws.Range("F47").Value = 1290
$ws.Range("F49").Value = 526
